$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph ("Play Blood And Gold for Free - A Thrilling
#    Medieval Game"). The new paragraph has a leading empty run (to
#    match the document's existing style), a bold "Meta description"
#    run, and a plain run with the rest of the description text.
# ----------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Blood And Gold for Free - A Thrilling Medieval Game</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Check out our review of Blood And Gold and play this exciting medieval game for free. Win big with bonuses and special features!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $titleRange.InsertXML($xml)

# ----------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicated bold paragraph
#    "Play Blood And Gold for Free - A Thrilling Medieval Game" (the
#    meta description was moved to the top of the document instead).
#    NOTE: the Heading1 title paragraph at the very top has the same
#    text, so search from the end of the document backwards and only
#    delete the plain "Normal"-style duplicate.
# ----------------------------------------------------------------------

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if (($para.Range.Text -eq "Play Blood And Gold for Free - A Thrilling Medieval Game`r") -and ($para.Style.NameLocal -eq "Normal")) {
        $para.Range.Delete()
        break
    }
}

# ----------------------------------------------------------------------
# 3) Replace the text of the remaining italic paragraph (previously the
#    meta description, now repurposed) with the DALLE image prompt,
#    keeping its italic formatting and the paragraph's leading empty
#    run intact.
# ----------------------------------------------------------------------

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Check out our review of Blood And Gold and play this exciting medieval game for free. Win big with bonuses and special features!`r") {
        $s = $para.Range.Start
        $e = $para.Range.End
        $textRange = $d.Range($s, $e - 1)
        $textRange.Text = 'Prompt for DALLE: Create a cartoon-style image for the game "Blood and Gold" that features a happy Maya warrior with glasses. The image should be eye-catching and dynamic, with vibrant colors that represent the game''s medieval theme. The Maya warrior should be holding a golden shield and a sword, with a determined expression on their face. The background should showcase a clash between two kingdoms, with catapults and horn sounds in the distance. Make sure the image conveys the game''s exciting theme and encourages players to try it out.'
        break
    }
}
